# edit.ps1 - applies the resume update described by the commit:
# "Updated Resume and Cover Letter (Eng & Kor Ver)"
#
# Changes:
#  1. Title "Chaewan Woo": drop the w:w=125 character-scaling and switch
#     to an explicit sz/szCs=54 (27pt) font size, split across two runs.
#  2. Refresh the w14:anchorId GUIDs on the four horizontal-rule VML
#     fallback shapes (Education / Skills / Projects / Competitive Gaming
#     headings) - a cosmetic re-save artifact, applied via a full-paragraph
#     InsertXML since that markup lives in the mc:Fallback and is not part
#     of the editable text story.
#  3. Move the stray "_GoBack" bookmark from inside the "Relevant
#     Coursework" run (where it split "Artificial Intelligence" into two
#     runs) to right after "Ontario College Advanced Diploma: Game
#     Programming", and re-join the coursework text into a single run.

$d = $word.ActiveDocument

# --- 1. Title: "Chaewan Woo" -> two runs, sz/szCs=54, no w:w ---
$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="008E7B35" w:rsidRPr="00912893" w:rsidRDefault="00B060CD" w:rsidP="00912893"><w:pPr><w:pStyle w:val="Title"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="54"/><w:szCs w:val="54"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="54"/><w:szCs w:val="54"/></w:rPr><w:t>Chaewan W</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="54"/><w:szCs w:val="54"/></w:rPr><w:t>oo</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(1).Range.InsertXML($titleXml)

# --- 3a. "Ontario College Advanced Diploma..." paragraph: add the
#         _GoBack bookmark right after the diploma-name run ---
$ontarioXml = @'
<w:p w:rsidR="00D01B20" w:rsidRPr="00F32272" w:rsidRDefault="00D01B20" w:rsidP="00D01B20"><w:pPr><w:tabs><w:tab w:val="left" w:pos="8664"/></w:tabs><w:spacing w:before="21"/><w:ind w:left="216"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr></w:pPr><w:r w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t>Ontario College Advanced Diploma: Game Programming</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t>Sep</w:t></w:r><w:r w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:spacing w:val="-4"/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t>2021</w:t></w:r><w:r w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:spacing w:val="18"/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="바탕" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t>–</w:t></w:r><w:r w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:spacing w:val="17"/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t>May</w:t></w:r><w:r w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:spacing w:val="-4"/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t>. 2024</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(6).Range.InsertXML($ontarioXml)

# --- 3b. "Relevant Coursework:" paragraph: remove the old _GoBack
#         bookmark and re-join the split "Artificial Intelligence" run ---
$courseworkXml = @'
<w:p w:rsidR="00D01B20" w:rsidRPr="00F32272" w:rsidRDefault="000E0CB8" w:rsidP="00D01B20"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="696"/></w:tabs><w:spacing w:line="235" w:lineRule="auto"/><w:ind w:right="559"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:position w:val="2"/><w:sz w:val="12"/><w:lang w:eastAsia="ko-KR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="바탕" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t xml:space="preserve">Relevant Coursework: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:w w:val="105"/><w:sz w:val="20"/></w:rPr><w:t>Data Structures &amp; Algorithms, Artificial Intelligence, Game Physics, Game Engines, 3D Graphics &amp; Rendering, Multiplayer Systems, Mobile &amp; Console Development</w:t></w:r><w:r w:rsidR="00D01B20" w:rsidRPr="00F32272"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="PMingLiU-ExtB" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:w w:val="105"/><w:sz w:val="20"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(9).Range.InsertXML($courseworkXml)

# --- 2. Refresh the four horizontal-rule shapes' w14:anchorId ---
$educationXml = @'
<w:p w:rsidR="009E620A" w:rsidRPr="008B432E" w:rsidRDefault="00E17B3D"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="213"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r w:rsidRPr="008B432E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:lang w:eastAsia="ko-KR"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="0" distR="0" simplePos="0" relativeHeight="251656192" behindDoc="1" locked="0" layoutInCell="1" allowOverlap="1"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="page"><wp:posOffset>457200</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>335983</wp:posOffset></wp:positionV><wp:extent cx="6858000" cy="1270"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapTopAndBottom/><wp:docPr id="4" name="Graphic 4"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr><a:spLocks/></wps:cNvSpPr><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="6858000" cy="1270"/></a:xfrm><a:custGeom><a:avLst/><a:gdLst/><a:ahLst/><a:cxnLst/><a:rect l="l" t="t" r="r" b="b"/><a:pathLst><a:path w="6858000"><a:moveTo><a:pt x="0" y="0"/></a:moveTo><a:lnTo><a:pt x="6858000" y="0"/></a:lnTo></a:path></a:pathLst></a:custGeom><a:ln w="5054"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:prstDash val="solid"/></a:ln></wps:spPr><wps:bodyPr wrap="square" lIns="0" tIns="0" rIns="0" bIns="0" rtlCol="0"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="71474945" id="Graphic 4" o:spid="_x0000_s1026" style="position:absolute;left:0;text-align:left;margin-left:36pt;margin-top:26.45pt;width:540pt;height:.1pt;z-index:-251660288;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:0;mso-wrap-distance-top:0;mso-wrap-distance-right:0;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:page;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:top" coordsize="6858000,1270" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQBK5LLYIQIAAH8EAAAOAAAAZHJzL2Uyb0RvYy54bWysVMFu2zAMvQ/YPwi6L3aCpAuMOMXQoMWA&#xA;oivQFDsrshwbk0WNUmLn70fJdpJ2t6I+CJT4RPLxUV7ddo1mR4WuBpPz6STlTBkJRW32OX/d3n9b&#xA;cua8MIXQYFTOT8rx2/XXL6vWZmoGFehCIaMgxmWtzXnlvc2SxMlKNcJNwCpDzhKwEZ62uE8KFC1F&#xA;b3QyS9ObpAUsLIJUztHppnfydYxflkr6X2XplGc651SbjyvGdRfWZL0S2R6FrWo5lCE+UEUjakNJ&#xA;z6E2wgt2wPq/UE0tERyUfiKhSaAsa6kiB2IzTd+xeamEVZELNcfZc5vc54WVT8dnZHWR8zlnRjQk&#xA;0cPQjXloTmtdRpgX+4yBnrOPIP84ciRvPGHjBkxXYhOwRI51sdOnc6dV55mkw5vlYpmmJIgk33T2&#xA;PQqRiGy8Kw/OPyiIccTx0flep2K0RDVasjOjiaR20FlHnT1npDNyRjrvep2t8OFeKC6YrL0UEs4a&#xA;OKotRK9/VzmVdvFqc406UxlZErZHkBHSUK96I6Ym+5qcNqGKRbqYx/FxoOvivtY6VOFwv7vTyI4i&#xA;DG/8Ag+K8AZm0fmNcFWPi64Bps2gUy9NEGkHxYkEb0njnLu/B4GKM/3T0EiF5zEaOBq70UCv7yA+&#xA;otggyrntfgu0LKTPuSdln2AcWJGNogXqZ2y4aeDHwUNZB0XjDPUVDRua8khweJHhGV3vI+ry31j/&#xA;AwAA//8DAFBLAwQUAAYACAAAACEAjGfI2N4AAAAJAQAADwAAAGRycy9kb3ducmV2LnhtbEyPQUvE&#xA;MBCF74L/IYzgRdy0la5ubbqIoKh4sSvoMduMbbGZlCS7rf/e6UmP897jm/fK7WwHcUQfekcK0lUC&#xA;AqlxpqdWwfvu4fIGRIiajB4coYIfDLCtTk9KXRg30Rse69gKhlAotIIuxrGQMjQdWh1WbkRi78t5&#xA;qyOfvpXG64nhdpBZkqyl1T3xh06PeN9h810fLFPWH+347KfNxdPuM3+psUns46tS52fz3S2IiHP8&#xA;C8NSn6tDxZ327kAmiEHBdcZTooI824BY/DRflD0rVynIqpT/F1S/AAAA//8DAFBLAQItABQABgAI&#xA;AAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsB&#xA;Ai0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsB&#xA;Ai0AFAAGAAgAAAAhAErkstghAgAAfwQAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1s&#xA;UEsBAi0AFAAGAAgAAAAhAIxnyNjeAAAACQEAAA8AAAAAAAAAAAAAAAAAewQAAGRycy9kb3ducmV2&#xA;LnhtbFBLBQYAAAAABAAEAPMAAACGBQAAAAA=&#xA;" path="m,l6858000,e" filled="f" strokeweight=".14039mm"><v:path arrowok="t"/><w10:wrap type="topAndBottom" anchorx="page"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:bookmarkStart w:id="0" w:name="Education"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="008B432E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:spacing w:val="-2"/><w:w w:val="130"/></w:rPr><w:t>E</w:t></w:r><w:r w:rsidRPr="008B432E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:smallCaps/><w:spacing w:val="-2"/><w:w w:val="130"/></w:rPr><w:t>ducation</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(4).Range.InsertXML($educationXml)

$skillsXml = @'
<w:p w:rsidR="009E620A" w:rsidRPr="008B432E" w:rsidRDefault="00E17B3D"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r w:rsidRPr="008B432E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:lang w:eastAsia="ko-KR"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="0" distR="0" simplePos="0" relativeHeight="251660288" behindDoc="1" locked="0" layoutInCell="1" allowOverlap="1"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="page"><wp:posOffset>457200</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>303980</wp:posOffset></wp:positionV><wp:extent cx="6858000" cy="1270"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapTopAndBottom/><wp:docPr id="5" name="Graphic 5"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr><a:spLocks/></wps:cNvSpPr><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="6858000" cy="1270"/></a:xfrm><a:custGeom><a:avLst/><a:gdLst/><a:ahLst/><a:cxnLst/><a:rect l="l" t="t" r="r" b="b"/><a:pathLst><a:path w="6858000"><a:moveTo><a:pt x="0" y="0"/></a:moveTo><a:lnTo><a:pt x="6858000" y="0"/></a:lnTo></a:path></a:pathLst></a:custGeom><a:ln w="5054"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:prstDash val="solid"/></a:ln></wps:spPr><wps:bodyPr wrap="square" lIns="0" tIns="0" rIns="0" bIns="0" rtlCol="0"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="15D3B657" id="Graphic 5" o:spid="_x0000_s1026" style="position:absolute;left:0;text-align:left;margin-left:36pt;margin-top:23.95pt;width:540pt;height:.1pt;z-index:-251656192;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:0;mso-wrap-distance-top:0;mso-wrap-distance-right:0;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:page;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:top" coordsize="6858000,1270" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQAtCWCrIQIAAH8EAAAOAAAAZHJzL2Uyb0RvYy54bWysVMFu2zAMvQ/YPwi6L3aCpSuMOMXQoMWA&#xA;oivQDDsrshwLk0WNUmLn70fJdpJ2t2E+CJT4RPLxUV7d9a1hR4Vegy35fJZzpqyEStt9yX9sHz7d&#xA;cuaDsJUwYFXJT8rzu/XHD6vOFWoBDZhKIaMg1hedK3kTgiuyzMtGtcLPwClLzhqwFYG2uM8qFB1F&#xA;b022yPObrAOsHIJU3tPpZnDydYpf10qG73XtVWCm5FRbSCumdRfXbL0SxR6Fa7QcyxD/UEUrtKWk&#xA;51AbEQQ7oP4rVKslgoc6zCS0GdS1lipxIDbz/B2b10Y4lbhQc7w7t8n/v7Dy+fiCTFclX3JmRUsS&#xA;PY7dWMbmdM4XhHl1LxjpefcE8pcnR/bGEzd+xPQ1thFL5FifOn06d1r1gUk6vLld3uY5CSLJN198&#xA;SUJkopjuyoMPjwpSHHF88mHQqZos0UyW7O1kIqkddTZJ58AZ6Yyckc67QWcnQrwXi4sm6y6FxLMW&#xA;jmoLyRveVU6lXbzGXqPOVCaWhB0QZMQ01KvBSKnJviZnbKximS8/p/HxYHT1oI2JVXjc7+4NsqOI&#xA;w5u+yIMivIE59GEjfDPgkmuEGTvqNEgTRdpBdSLBO9K45P73QaDizHyzNFLxeUwGTsZuMjCYe0iP&#xA;KDWIcm77nwIdi+lLHkjZZ5gGVhSTaJH6GRtvWvh6CFDrqGiaoaGicUNTngiOLzI+o+t9Ql3+G+s/&#xA;AAAA//8DAFBLAwQUAAYACAAAACEAanTR7d4AAAAJAQAADwAAAGRycy9kb3ducmV2LnhtbEyPzU7D&#xA;MBCE70i8g7VIXBB1UtG/EKdCSCBAXEiR2qMbL0lEvI5stwlvz+YEx50ZfTuTb0fbiTP60DpSkM4S&#xA;EEiVMy3VCj53T7drECFqMrpzhAp+MMC2uLzIdWbcQB94LmMtGEIh0wqaGPtMylA1aHWYuR6JvS/n&#xA;rY58+loarweG207Ok2QprW6JPzS6x8cGq+/yZJmy3Nf9qx82Ny+7w+KtxCqxz+9KXV+ND/cgIo7x&#xA;LwxTfa4OBXc6uhOZIDoFqzlPiQruVhsQk58uJuXIyjoFWeTy/4LiFwAA//8DAFBLAQItABQABgAI&#xA;AAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsB&#xA;Ai0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsB&#xA;Ai0AFAAGAAgAAAAhAC0JYKshAgAAfwQAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1s&#xA;UEsBAi0AFAAGAAgAAAAhAGp00e3eAAAACQEAAA8AAAAAAAAAAAAAAAAAewQAAGRycy9kb3ducmV2&#xA;LnhtbFBLBQYAAAAABAAEAPMAAACGBQAAAAA=&#xA;" path="m,l6858000,e" filled="f" strokeweight=".14039mm"><v:path arrowok="t"/><w10:wrap type="topAndBottom" anchorx="page"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:bookmarkStart w:id="2" w:name="Skills"/><w:bookmarkEnd w:id="2"/><w:r w:rsidRPr="008B432E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:spacing w:val="-2"/><w:w w:val="145"/></w:rPr><w:t>S</w:t></w:r><w:r w:rsidRPr="008B432E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:smallCaps/><w:spacing w:val="-2"/><w:w w:val="145"/></w:rPr><w:t>kills</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(18).Range.InsertXML($skillsXml)

$projectsXml = @'
<w:p w:rsidR="009E620A" w:rsidRPr="008B432E" w:rsidRDefault="00E17B3D"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r w:rsidRPr="008B432E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:lang w:eastAsia="ko-KR"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="0" distR="0" simplePos="0" relativeHeight="251669504" behindDoc="1" locked="0" layoutInCell="1" allowOverlap="1"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="page"><wp:posOffset>457200</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>303709</wp:posOffset></wp:positionV><wp:extent cx="6858000" cy="1270"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapTopAndBottom/><wp:docPr id="7" name="Graphic 7"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr><a:spLocks/></wps:cNvSpPr><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="6858000" cy="1270"/></a:xfrm><a:custGeom><a:avLst/><a:gdLst/><a:ahLst/><a:cxnLst/><a:rect l="l" t="t" r="r" b="b"/><a:pathLst><a:path w="6858000"><a:moveTo><a:pt x="0" y="0"/></a:moveTo><a:lnTo><a:pt x="6858000" y="0"/></a:lnTo></a:path></a:pathLst></a:custGeom><a:ln w="5054"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:prstDash val="solid"/></a:ln></wps:spPr><wps:bodyPr wrap="square" lIns="0" tIns="0" rIns="0" bIns="0" rtlCol="0"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="35B32C5F" id="Graphic 7" o:spid="_x0000_s1026" style="position:absolute;left:0;text-align:left;margin-left:36pt;margin-top:23.9pt;width:540pt;height:.1pt;z-index:-251646976;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:0;mso-wrap-distance-top:0;mso-wrap-distance-right:0;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:page;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:top" coordsize="6858000,1270" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQDj08VMIQIAAH8EAAAOAAAAZHJzL2Uyb0RvYy54bWysVMFu2zAMvQ/YPwi6L3aCpQmMOMXQoMWA&#xA;oivQFDsrshwbk0WNUmLn70fJdpJ2t6I+CJT4RPLxUV7ddo1mR4WuBpPz6STlTBkJRW32OX/d3n9b&#xA;cua8MIXQYFTOT8rx2/XXL6vWZmoGFehCIaMgxmWtzXnlvc2SxMlKNcJNwCpDzhKwEZ62uE8KFC1F&#xA;b3QyS9ObpAUsLIJUztHppnfydYxflkr6X2XplGc651SbjyvGdRfWZL0S2R6FrWo5lCE+UEUjakNJ&#xA;z6E2wgt2wPq/UE0tERyUfiKhSaAsa6kiB2IzTd+xeamEVZELNcfZc5vc54WVT8dnZHWR8wVnRjQk&#xA;0cPQjUVoTmtdRpgX+4yBnrOPIP84ciRvPGHjBkxXYhOwRI51sdOnc6dV55mkw5vlfJmmJIgk33S2&#xA;iEIkIhvvyoPzDwpiHHF8dL7XqRgtUY2W7MxoIqkddNZRZ88Z6Yyckc67XmcrfLgXigsmay+FhLMG&#xA;jmoL0evfVU6lXbzaXKPOVEaWhO0RZIQ01KveiKnJvianTahins6/x/FxoOvivtY6VOFwv7vTyI4i&#xA;DG/8Ag+K8AZm0fmNcFWPi64Bps2gUy9NEGkHxYkEb0njnLu/B4GKM/3T0EiF5zEaOBq70UCv7yA+&#xA;otggyrntfgu0LKTPuSdln2AcWJGNogXqZ2y4aeDHwUNZB0XjDPUVDRua8khweJHhGV3vI+ry31j/&#xA;AwAA//8DAFBLAwQUAAYACAAAACEALePDut4AAAAJAQAADwAAAGRycy9kb3ducmV2LnhtbEyPzU7D&#xA;MBCE70i8g7VIXBC1W9EfQpwKIYGg4kKKBEc3XpKIeB3ZbhPens0Jjjszmp0v346uEycMsfWkYT5T&#xA;IJAqb1uqNbzvH683IGIyZE3nCTX8YIRtcX6Wm8z6gd7wVKZacAnFzGhoUuozKWPVoDNx5nsk9r58&#xA;cCbxGWppgxm43HVyodRKOtMSf2hMjw8NVt/l0XHL6qPuX8Jwe/W8/1zuSqyUe3rV+vJivL8DkXBM&#xA;f2GY5vN0KHjTwR/JRtFpWC8YJWm4WTPB5M+Xk3JgZaNAFrn8T1D8AgAA//8DAFBLAQItABQABgAI&#xA;AAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsB&#xA;Ai0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsB&#xA;Ai0AFAAGAAgAAAAhAOPTxUwhAgAAfwQAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1s&#xA;UEsBAi0AFAAGAAgAAAAhAC3jw7reAAAACQEAAA8AAAAAAAAAAAAAAAAAewQAAGRycy9kb3ducmV2&#xA;LnhtbFBLBQYAAAAABAAEAPMAAACGBQAAAAA=&#xA;" path="m,l6858000,e" filled="f" strokeweight=".14039mm"><v:path arrowok="t"/><w10:wrap type="topAndBottom" anchorx="page"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:bookmarkStart w:id="3" w:name="Projects"/><w:bookmarkEnd w:id="3"/><w:r w:rsidR="00EB302D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:smallCaps/><w:spacing w:val="-2"/><w:w w:val="140"/></w:rPr><w:t>Projects</w:t></w:r><w:r w:rsidR="00D1432B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:smallCaps/><w:spacing w:val="-2"/><w:w w:val="140"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00693F66"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:smallCaps/><w:spacing w:val="-2"/><w:w w:val="140"/></w:rPr><w:t xml:space="preserve">| </w:t></w:r><w:r w:rsidR="00CA10B3"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:w w:val="115"/></w:rPr><w:t>https://w</w:t></w:r><w:r w:rsidR="0047020A" w:rsidRPr="0047020A"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:w w:val="115"/></w:rPr><w:t>oo95.github.io/</w:t></w:r><w:r w:rsidR="00696DBC"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:w w:val="115"/></w:rPr><w:t>Devhub</w:t></w:r><w:r w:rsidR="00CA10B3"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:w w:val="115"/></w:rPr><w:t>/</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(24).Range.InsertXML($projectsXml)

$competitiveXml = @'
<w:p w:rsidR="00D71ED0" w:rsidRDefault="00EB302D" w:rsidP="00D71ED0"><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="229"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:smallCaps/><w:spacing w:val="-2"/><w:w w:val="140"/></w:rPr></w:pPr><w:r w:rsidRPr="008B432E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:lang w:eastAsia="ko-KR"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="0" distR="0" simplePos="0" relativeHeight="251676672" behindDoc="1" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="5222E330" wp14:editId="1F729563"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="page"><wp:posOffset>457200</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>346290</wp:posOffset></wp:positionV><wp:extent cx="6858000" cy="1270"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapTopAndBottom/><wp:docPr id="9" name="Graphic 6"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr><a:spLocks/></wps:cNvSpPr><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="6858000" cy="1270"/></a:xfrm><a:custGeom><a:avLst/><a:gdLst/><a:ahLst/><a:cxnLst/><a:rect l="l" t="t" r="r" b="b"/><a:pathLst><a:path w="6858000"><a:moveTo><a:pt x="0" y="0"/></a:moveTo><a:lnTo><a:pt x="6858000" y="0"/></a:lnTo></a:path></a:pathLst></a:custGeom><a:ln w="5054"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:prstDash val="solid"/></a:ln></wps:spPr><wps:bodyPr wrap="square" lIns="0" tIns="0" rIns="0" bIns="0" rtlCol="0"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="45AC6B8F" id="Graphic 6" o:spid="_x0000_s1026" style="position:absolute;left:0;text-align:left;margin-left:36pt;margin-top:27.25pt;width:540pt;height:.1pt;z-index:-251639808;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:0;mso-wrap-distance-top:0;mso-wrap-distance-right:0;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:page;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:top" coordsize="6858000,1270" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQA6Bj7RIQIAAH8EAAAOAAAAZHJzL2Uyb0RvYy54bWysVMFu2zAMvQ/YPwi6L3aCJUuNOMXQoMWA&#xA;oivQDDsrshwblSWNVGLn70fJdpJ2t6E+CJT4RPLxUV7ddo1mRwVYW5Pz6STlTBlpi9rsc/5re/9l&#xA;yRl6YQqhrVE5Pynkt+vPn1aty9TMVlYXChgFMZi1LueV9y5LEpSVagROrFOGnKWFRnjawj4pQLQU&#xA;vdHJLE0XSWuhcGClQqTTTe/k6xi/LJX0P8sSlWc651SbjyvEdRfWZL0S2R6Eq2o5lCH+o4pG1IaS&#xA;nkNthBfsAPU/oZpagkVb+om0TWLLspYqciA20/Qdm5dKOBW5UHPQnduEHxdWPh2fgdVFzm84M6Ih&#xA;iR6GbixCc1qHGWFe3DMEeugerXxFciRvPGGDA6YroQlYIse62OnTudOq80zS4WI5X6YpCSLJN519&#xA;i0IkIhvvygP6B2VjHHF8RN/rVIyWqEZLdmY0gdQOOuuos+eMdAbOSOddr7MTPtwLxQWTtZdCwllj&#xA;j2pro9e/q5xKu3i1uUadqYwsCdsjyAhpqFe9EVOTfU1Om1DFPJ1/jeODVtfFfa11qAJhv7vTwI4i&#xA;DG/8Ag+K8AbmAP1GYNXjomuAaTPo1EsTRNrZ4kSCt6RxzvHPQYDiTP8wNFLheYwGjMZuNMDrOxsf&#xA;UWwQ5dx2vwU4FtLn3JOyT3YcWJGNogXqZ2y4aez3g7dlHRSNM9RXNGxoyiPB4UWGZ3S9j6jLf2P9&#xA;FwAA//8DAFBLAwQUAAYACAAAACEAAzF1jd4AAAAJAQAADwAAAGRycy9kb3ducmV2LnhtbEyPwU7D&#xA;MBBE70j9B2srcUHUaUVaGuJUCAkEiAspEj268ZJExOvIdpvw92xOcNyZ0duZfDfaTpzRh9aRguUi&#xA;AYFUOdNSreBj/3h9CyJETUZ3jlDBDwbYFbOLXGfGDfSO5zLWgiEUMq2gibHPpAxVg1aHheuR2Pty&#xA;3urIp6+l8XpguO3kKknW0uqW+EOje3xosPouT5Yp68+6f/HD9up5f0hfS6wS+/Sm1OV8vL8DEXGM&#xA;f2GY6nN1KLjT0Z3IBNEp2Kx4SlSQ3qQgJn+ZTspxUjYgi1z+X1D8AgAA//8DAFBLAQItABQABgAI&#xA;AAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsB&#xA;Ai0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsB&#xA;Ai0AFAAGAAgAAAAhADoGPtEhAgAAfwQAAA4AAAAAAAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1s&#xA;UEsBAi0AFAAGAAgAAAAhAAMxdY3eAAAACQEAAA8AAAAAAAAAAAAAAAAAewQAAGRycy9kb3ducmV2&#xA;LnhtbFBLBQYAAAAABAAEAPMAAACGBQAAAAA=&#xA;" path="m,l6858000,e" filled="f" strokeweight=".14039mm"><v:path arrowok="t"/><w10:wrap type="topAndBottom" anchorx="page"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:smallCaps/><w:spacing w:val="-2"/><w:w w:val="140"/></w:rPr><w:t>Competitive Gaming Achievements</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(38).Range.InsertXML($competitiveXml)
